$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns to reflect that the year fields represent the
#    season-ending calendar year rather than a generic "year".
$ws.Range("A1").Value = "season_ending_year_x"
$ws.Range("O1").Value = "season_ending_year_y"

# 2) Add a new trailing column AY: "calendar_year" which mirrors the
#    season_ending_year_x (column A) value for each data row.
$ws.Range("AY1").Value = "calendar_year"
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)

$calendarYears = @(1976, 1975, 1975, 1974, 1973, 1972, 1971, 1970, 1969, 1968)
for ($i = 0; $i -lt $calendarYears.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 51).Value = $calendarYears[$i]
}

# 3) Update birth_year (column Q) values: increment each by one year.
$birthYears = @(1951, 1951, 1951, 1951, 1944, 1950, 1945, 1950, 1945, 1943)
for ($i = 0; $i -lt $birthYears.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $birthYears[$i]
}

Write-Output "done"
